$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on humidity (%) cells so literal values like "71%" are preserved as text
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H46").NumberFormat = "@"

$ws.Range("E2").Value = "2026-02-15 14:48:24"
$ws.Range("E3").Value = "2026-02-15 14:48:26"
$ws.Range("I3").Value = "0.8 mm"
$ws.Range("K3").Value = "6.0 MJ/m2"
$ws.Range("O3").Value = "-6.3 °C"
$ws.Range("E4").Value = "2026-02-15 14:48:29"
$ws.Range("H4").Value = "71%"
$ws.Range("K4").Value = "10.4 MJ/m2"
$ws.Range("O4").Value = "6.0 °C"
$ws.Range("E5").Value = "2026-02-15 14:48:31"
$ws.Range("H5").Value = "92%"
$ws.Range("I5").Value = "3.0 mm"
$ws.Range("K5").Value = "4.7 MJ/m2"
$ws.Range("O5").Value = "-5.7 °C"
$ws.Range("E6").Value = "2026-02-15 14:48:34"
$ws.Range("H6").Value = "58%"
$ws.Range("K6").Value = "10.9 MJ/m2"
$ws.Range("O6").Value = "7.5 °C"
$ws.Range("E7").Value = "2026-02-15 14:48:36"
$ws.Range("K7").Value = "10.9 MJ/m2"
$ws.Range("O7").Value = "11.1 °C"
$ws.Range("E8").Value = "2026-02-15 14:48:39"
$ws.Range("K8").Value = "10.9 MJ/m2"
$ws.Range("O8").Value = "7.5 °C"
$ws.Range("E9").Value = "2026-02-15 14:48:41"
$ws.Range("H9").Value = "44%"
$ws.Range("K9").Value = "10.7 MJ/m2"
$ws.Range("O9").Value = "10.8 °C"
$ws.Range("E10").Value = "2026-02-15 14:48:43"
$ws.Range("K10").Value = "11.0 MJ/m2"
$ws.Range("M10").Value = "13.9 °C 14:10 TU"
$ws.Range("O10").Value = "6.6 °C"
$ws.Range("E11").Value = "2026-02-15 14:48:46"
$ws.Range("H11").Value = "35%"
$ws.Range("O11").Value = "7.2 °C"
$ws.Range("E12").Value = "2026-02-15 14:48:48"
$ws.Range("O12").Value = "10.6 °C"
$ws.Range("E13").Value = "2026-02-15 14:48:50"
$ws.Range("J13").Value = "1015.9 hPa"
$ws.Range("K13").Value = "5.9 MJ/m2"
$ws.Range("M13").Value = "12.5 °C 14:02 TU"
$ws.Range("O13").Value = "5.6 °C"
$ws.Range("E14").Value = "2026-02-15 14:48:53"
$ws.Range("K14").Value = "10.5 MJ/m2"
$ws.Range("O14").Value = "10.5 °C"
$ws.Range("E15").Value = "2026-02-15 14:48:55"
$ws.Range("O15").Value = "10.5 °C"
$ws.Range("E16").Value = "2026-02-15 14:48:58"
$ws.Range("H16").Value = "57%"
$ws.Range("K16").Value = "8.3 MJ/m2"
$ws.Range("O16").Value = "-2.6 °C"
$ws.Range("E17").Value = "2026-02-15 14:49:00"
$ws.Range("H17").Value = "28%"
$ws.Range("K17").Value = "11.6 MJ/m2"
$ws.Range("O17").Value = "3.0 °C"
$ws.Range("E18").Value = "2026-02-15 14:49:03"
$ws.Range("H18").Value = "75%"
$ws.Range("K18").Value = "11.0 MJ/m2"
$ws.Range("O18").Value = "5.9 °C"
$ws.Range("E19").Value = "2026-02-15 14:49:05"
$ws.Range("K19").Value = "10.8 MJ/m2"
$ws.Range("M19").Value = "8.9 °C 14:01 TU"
$ws.Range("O19").Value = "2.2 °C"
$ws.Range("E20").Value = "2026-02-15 14:49:07"
$ws.Range("H20").Value = "53%"
$ws.Range("K20").Value = "11.7 MJ/m2"
$ws.Range("O20").Value = "-3.5 °C"
$ws.Range("E21").Value = "2026-02-15 14:49:10"
$ws.Range("J21").Value = "1015.2 hPa"
$ws.Range("K21").Value = "9.7 MJ/m2"
$ws.Range("M21").Value = "15.0 °C 14:01 TU"
$ws.Range("O21").Value = "7.1 °C"
$ws.Range("E22").Value = "2026-02-15 14:49:13"
$ws.Range("K22").Value = "11.3 MJ/m2"
$ws.Range("N22").Value = "-6.5 °C 14:00 TU"
$ws.Range("E23").Value = "2026-02-15 14:49:15"
$ws.Range("H23").Value = "59%"
$ws.Range("K23").Value = "11.6 MJ/m2"
$ws.Range("O23").Value = "-4.5 °C"
$ws.Range("E24").Value = "2026-02-15 14:49:18"
$ws.Range("K24").Value = "10.3 MJ/m2"
$ws.Range("L24").Value = "51.5 km/h - 296º 14:29 TU"
$ws.Range("O24").Value = "7.7 °C"
$ws.Range("E25").Value = "2026-02-15 14:49:20"
$ws.Range("H25").Value = "58%"
$ws.Range("K25").Value = "8.1 MJ/m2"
$ws.Range("M25").Value = "2.2 °C 14:03 TU"
$ws.Range("O25").Value = "-2.5 °C"
$ws.Range("E26").Value = "2026-02-15 14:49:23"
$ws.Range("E27").Value = "2026-02-15 14:49:25"
$ws.Range("H27").Value = "41%"
$ws.Range("K27").Value = "10.1 MJ/m2"
$ws.Range("M27").Value = "3.1 °C 14:13 TU"
$ws.Range("O27").Value = "-0.7 °C"
$ws.Range("E28").Value = "2026-02-15 14:49:27"
$ws.Range("K28").Value = "9.7 MJ/m2"
$ws.Range("L28").Value = "16.6 km/h - 26º 14:05 TU"
$ws.Range("O28").Value = "5.4 °C"
$ws.Range("E29").Value = "2026-02-15 14:49:29"
$ws.Range("K29").Value = "11.3 MJ/m2"
$ws.Range("O29").Value = "9.9 °C"
$ws.Range("E30").Value = "2026-02-15 14:49:32"
$ws.Range("K30").Value = "11.0 MJ/m2"
$ws.Range("O30").Value = "9.5 °C"
$ws.Range("E31").Value = "2026-02-15 14:49:34"
$ws.Range("J31").Value = "1014.1 hPa"
$ws.Range("K31").Value = "9.5 MJ/m2"
$ws.Range("O31").Value = "8.9 °C"
$ws.Range("E32").Value = "2026-02-15 14:49:36"
$ws.Range("H32").Value = "87%"
$ws.Range("K32").Value = "8.2 MJ/m2"
$ws.Range("O32").Value = "2.5 °C"
$ws.Range("E33").Value = "2026-02-15 14:49:39"
$ws.Range("J33").Value = "1015.8 hPa"
$ws.Range("K33").Value = "9.7 MJ/m2"
$ws.Range("O33").Value = "4.9 °C"
$ws.Range("E34").Value = "2026-02-15 14:49:41"
$ws.Range("K34").Value = "11.1 MJ/m2"
$ws.Range("O34").Value = "0.1 °C"
$ws.Range("E35").Value = "2026-02-15 14:49:44"
$ws.Range("J35").Value = "1019.7 hPa"
$ws.Range("K35").Value = "8.5 MJ/m2"
$ws.Range("O35").Value = "3.1 °C"
$ws.Range("E36").Value = "2026-02-15 14:49:46"
$ws.Range("K36").Value = "9.0 MJ/m2"
$ws.Range("O36").Value = "11.0 °C"
$ws.Range("E37").Value = "2026-02-15 14:49:49"
$ws.Range("H37").Value = "50%"
$ws.Range("O37").Value = "5.4 °C"
$ws.Range("E38").Value = "2026-02-15 14:49:51"
$ws.Range("H38").Value = "66%"
$ws.Range("K38").Value = "11.2 MJ/m2"
$ws.Range("O38").Value = "6.6 °C"
$ws.Range("E39").Value = "2026-02-15 14:49:54"
$ws.Range("K39").Value = "8.7 MJ/m2"
$ws.Range("O39").Value = "-4.0 °C"
$ws.Range("E40").Value = "2026-02-15 14:49:56"
$ws.Range("O40").Value = "8.6 °C"
$ws.Range("E41").Value = "2026-02-15 14:49:59"
$ws.Range("H41").Value = "50%"
$ws.Range("K41").Value = "11.2 MJ/m2"
$ws.Range("O41").Value = "11.4 °C"
$ws.Range("E42").Value = "2026-02-15 14:50:01"
$ws.Range("O42").Value = "10.3 °C"
$ws.Range("E43").Value = "2026-02-15 14:50:04"
$ws.Range("H43").Value = "69%"
$ws.Range("K43").Value = "11.6 MJ/m2"
$ws.Range("O43").Value = "4.9 °C"
$ws.Range("E44").Value = "2026-02-15 14:50:06"
$ws.Range("K44").Value = "8.5 MJ/m2"
$ws.Range("O44").Value = "-5.0 °C"
$ws.Range("E45").Value = "2026-02-15 14:50:09"
$ws.Range("I45").Value = "0.2 mm"
$ws.Range("J45").Value = "1024.5 hPa"
$ws.Range("K45").Value = "3.7 MJ/m2"
$ws.Range("O45").Value = "-0.1 °C"
$ws.Range("E46").Value = "2026-02-15 14:50:11"
$ws.Range("H46").Value = "53%"
$ws.Range("J46").Value = "1019.9 hPa"
$ws.Range("K46").Value = "10.6 MJ/m2"
$ws.Range("O46").Value = "10.7 °C"
